# Applies the scheduled-runner data refresh to the Leve profit tables.
# For each sheet, numeric market-data cells (H:N) are updated to the latest
# scraped values. Where the refreshed profit figure is blank/not-applicable,
# the cell is cleared entirely (matching how this workbook already represents
# missing NQ/HQ profit figures elsewhere in the sheet), and conversely a cell
# that now has a computable figure is written in.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 551.8
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 551.8
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1655.4
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = -1991.4
$ws.Range("H98").Value = 947.44446
$ws.Range("I98").Value = 947.24
$ws.Range("K98").Value = 947.24
$ws.Range("M98").Value = 550.76
$ws.Range("H122").Value = 947.44446
$ws.Range("I122").Value = 947.24
$ws.Range("K122").Value = 2841.72
$ws.Range("M122").Value = -391.7200000000003
$ws.Range("H137").Value = 6417866.5
$ws.Range("I137").Value = 11908347
$ws.Range("K137").Value = 35725041
$ws.Range("M137").Value = -35722491

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6136440.5
$ws.Range("I2").Value = 11501786
$ws.Range("J2").Value = 4616.4287
$ws.Range("K2").Value = 11501786
$ws.Range("L2").Value = 4616.4287
$ws.Range("M2").Value = -11501673
$ws.Range("N2").Value = -4842.4287
$ws.Range("H32").Value = 32558.576
$ws.Range("I32").Value = 33497.906
$ws.Range("K32").Value = 33497.906
$ws.Range("M32").Value = -33210.906
$ws.Range("H61").Value = 3512.0476
$ws.Range("I61").Value = 2829.1052
$ws.Range("K61").Value = 2829.1052
$ws.Range("M61").Value = -2617.1052
$ws.Range("H116").Value = 6136440.5
$ws.Range("I116").Value = 11501786
$ws.Range("J116").Value = 4616.4287
$ws.Range("K116").Value = 11501786
$ws.Range("L116").Value = 4616.4287
$ws.Range("M116").Value = -11499492
$ws.Range("N116").Value = -9204.4287
$ws.Range("H132").Value = 3822.8333
$ws.Range("I132").Value = 2644.4375
$ws.Range("K132").Value = 7933.3125
$ws.Range("M132").Value = -5403.3125
$ws.Range("H136").Value = 3512.0476
$ws.Range("I136").Value = 2829.1052
$ws.Range("K136").Value = 8487.3156
$ws.Range("M136").Value = -5937.3156

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6136440.5
$ws.Range("I3").Value = 11501786
$ws.Range("J3").Value = 4616.4287
$ws.Range("K3").Value = 11501786
$ws.Range("L3").Value = 4616.4287
$ws.Range("M3").Value = -11501672
$ws.Range("N3").Value = -4844.4287
$ws.Range("H105").Value = 47633864
$ws.Range("I105").Value = 58840790
$ws.Range("J105").Value = 4425
$ws.Range("K105").Value = 58840790
$ws.Range("L105").Value = 4425
$ws.Range("M105").Value = -58839043
$ws.Range("N105").Value = -7919
$ws.Range("H107").Value = 3191.1667
$ws.Range("I107").Value = 3090.125
$ws.Range("J107").Value = 3999.5
$ws.Range("K107").Value = 3090.125
$ws.Range("L107").Value = 3999.5
$ws.Range("M107").Value = -1170.125
$ws.Range("N107").Value = -7839.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29414458
$ws.Range("I31").Value = 41667580
$ws.Range("J31").Value = 6964.9
$ws.Range("K31").Value = 41667580
$ws.Range("L31").Value = 6964.9
$ws.Range("M31").Value = -41667285
$ws.Range("N31").Value = -7554.9
$ws.Range("H34").Value = 29414458
$ws.Range("I34").Value = 41667580
$ws.Range("J34").Value = 6964.9
$ws.Range("K34").Value = 41667580
$ws.Range("L34").Value = 6964.9
$ws.Range("M34").Value = -41667378
$ws.Range("N34").Value = -7368.9
$ws.Range("H105").Value = 9973
$ws.Range("I105").Value = 20502.5
$ws.Range("K105").Value = 20502.5
$ws.Range("M105").Value = -18755.5
$ws.Range("H122").Value = 75861.64
$ws.Range("I122").Value = 81551
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 244653
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -242203
$ws.Range("N122").Value = -10600
$ws.Range("H134").Value = 6209.433
$ws.Range("I134").Value = 5263.52
$ws.Range("K134").Value = 15790.56
$ws.Range("M134").Value = -13255.56

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 4404.4375
$ws.Range("I2").Value = 63.909092
$ws.Range("J2").Value = 13953.6
$ws.Range("K2").Value = 383.454552
$ws.Range("L2").Value = 83721.60000000001
$ws.Range("M2").Value = -270.454552
$ws.Range("N2").Value = -83947.60000000001
$ws.Range("H7").Value = 80.833336
$ws.Range("I7").Value = 67.333336
$ws.Range("J7").Value = 94.333336
$ws.Range("K7").Value = 202.000008
$ws.Range("L7").Value = 283.000008
$ws.Range("M7").Value = -90.00000800000001
$ws.Range("N7").Value = -507.000008
$ws.Range("H32").Value = 700
$ws.Range("J32").Value = 700
$ws.Range("L32").Value = 2100
$ws.Range("N32").Value = -2666
$ws.Range("H34").Value = 1790.1428
$ws.Range("I34").Value = 70.57143000000001
$ws.Range("J34").Value = 2649.9285
$ws.Range("K34").Value = 211.71429
$ws.Range("L34").Value = 7949.7855
$ws.Range("M34").Value = -127.71429
$ws.Range("N34").Value = -8117.7855
$ws.Range("H92").Value = 2699.842
$ws.Range("I92").Value = 2766.6667
$ws.Range("J92").Value = 2639.7
$ws.Range("K92").Value = 8300.000100000001
$ws.Range("L92").Value = 7919.099999999999
$ws.Range("M92").Value = -7052.000100000001
$ws.Range("N92").Value = -10415.1
$ws.Range("H97").Value = 630.1429000000001
$ws.Range("J97").Value = 946.75
$ws.Range("L97").Value = 2840.25
$ws.Range("N97").Value = -3832.25
$ws.Range("H107").Value = 4663.143
$ws.Range("I107").Value = 550
$ws.Range("J107").Value = 5348.6665
$ws.Range("K107").Value = 1650
$ws.Range("L107").Value = 16045.9995
$ws.Range("M107").Value = 270
$ws.Range("N107").Value = -19885.9995
$ws.Range("H113").Value = 907.4
$ws.Range("J113").Value = 798
$ws.Range("L113").Value = 2394
$ws.Range("N113").Value = -6734
$ws.Range("H116").Value = 2009.8
$ws.Range("I116").Value = 2009.8
$ws.Range("K116").Value = 6029.4
$ws.Range("M116").Value = -2587.4
$ws.Range("H122").Value = 11964.167
$ws.Range("I122").Value = 23346.334
$ws.Range("J122").Value = 582
$ws.Range("K122").Value = 210117.006
$ws.Range("L122").Value = 5238
$ws.Range("M122").Value = -207667.006
$ws.Range("N122").Value = -10138
$ws.Range("H129").Value = 25003318
$ws.Range("I129").Value = 3492.5789
$ws.Range("J129").Value = 500000000
$ws.Range("K129").Value = 10477.7367
$ws.Range("L129").Value = 1500000000
$ws.Range("M129").Value = -5477.736699999999
$ws.Range("N129").Value = -1500010000
$ws.Range("H131").Value = 9808936
$ws.Range("I131").Value = 25641898
$ws.Range("J131").Value = 7579.1904
$ws.Range("K131").Value = 76925694
$ws.Range("L131").Value = 22737.5712
$ws.Range("M131").Value = -76920654
$ws.Range("N131").Value = -32817.57120000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4348.8667
$ws.Range("I80").Value = 4332.125
$ws.Range("K80").Value = 4332.125
$ws.Range("M80").Value = -3334.125
$ws.Range("H83").Value = 4348.8667
$ws.Range("I83").Value = 4332.125
$ws.Range("K83").Value = 21660.625
$ws.Range("M83").Value = -16668.625
$ws.Range("H97").Value = 1323.6111
$ws.Range("I97").Value = 1314.9231
$ws.Range("J97").Value = 1346.2
$ws.Range("K97").Value = 1314.9231
$ws.Range("L97").Value = 1346.2
$ws.Range("M97").Value = -818.9231
$ws.Range("N97").Value = -2338.2
$ws.Range("H122").Value = 11342.939
$ws.Range("I122").Value = 12462.27
$ws.Range("J122").Value = 7185.4287
$ws.Range("K122").Value = 37386.81
$ws.Range("L122").Value = 21556.2861
$ws.Range("M122").Value = -34936.81
$ws.Range("N122").Value = -26456.2861
$ws.Range("H133").Value = 176156
$ws.Range("J133").Value = 176156
$ws.Range("L133").Value = 176156
$ws.Range("N133").Value = -186276
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""
$ws.Range("H140").Value = 109997.5
$ws.Range("J140").Value = 109997.5
$ws.Range("L140").Value = 109997.5
$ws.Range("N140").Value = -120357.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 1006
$ws.Range("I9").Value = 1006
$ws.Range("K9").Value = 1006
$ws.Range("M9").Value = -866
$ws.Range("H113").Value = 1274.6666
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = ""
$ws.Range("H122").Value = 1910.2
$ws.Range("I122").Value = 1910.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5730.6
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3280.6
$ws.Range("N122").Value = ""
$ws.Range("H126").Value = 2332.9285
$ws.Range("I126").Value = 1412.96
$ws.Range("K126").Value = 4238.88
$ws.Range("M126").Value = -1768.88
$ws.Range("H136").Value = 2526.4
$ws.Range("I136").Value = 942.82355
$ws.Range("K136").Value = 2828.47065
$ws.Range("M136").Value = -278.4706499999998
